$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.58000000000056
$ws.Range("G2").Value = [double]"1.144397428931576e-06"
$ws.Range("H2").Value = [double]"2.013084047279802e-05"
$ws.Range("K2").Value = 5.555226070754964
$ws.Range("L2").Value = "[3.1038660376510805, 8.006586103858847]"
$ws.Range("M2").Value = [double]"1.100359731354139e-05"
$ws.Range("N2").Value = [double]"2.200719462708278e-05"
$ws.Range("O2").Value = -1.358526552903695
$ws.Range("P2").Value = "[-1.8742634850245423, -0.8427896207828471]"
$ws.Range("Q2").Value = [double]"3.616308816845049e-07"
$ws.Range("R2").Value = [double]"3.616308816845049e-07"
$ws.Range("S2").Value = 12.03585748987473
$ws.Range("T2").Value = "[10.588095873022208, 13.48361910672725]"
$ws.Range("W2").Value = 5.530810810810934
$ws.Range("X2").Value = 3.431151151151226
$ws.Range("Y2").Value = 7.630470470470642

# Row 3 updates
$ws.Range("E3").Value = 24.49000000000039
$ws.Range("G3").Value = 0.0007706462536779535
$ws.Range("H3").Value = 0.002802182823273261
$ws.Range("K3").Value = 4.480020723192665
$ws.Range("L3").Value = "[1.6743054995047473, 7.285735946880582]"
$ws.Range("M3").Value = 0.001849089378651358
$ws.Range("N3").Value = 0.001849089378651358
$ws.Range("O3").Value = 2.937184723176043
$ws.Range("P3").Value = "[2.157289850212811, 3.717079596139275]"
$ws.Range("Q3").Value = [double]"1.414868222582299e-12"
$ws.Range("R3").Value = [double]"2.829736445164599e-12"
$ws.Range("S3").Value = 11.31585656786661
$ws.Range("T3").Value = "[9.674590695016743, 12.957122440716468]"
$ws.Range("W3").Value = 13.04172172172193
$ws.Range("X3").Value = 10.00192192192208
$ws.Range("Y3").Value = 16.08152152152178
